# The data table (rows 2-49, columns A-H) is organised as four 12-month
# year blocks (2014, 2015, 2016, 2017), each block starting with January
# and ending with December. The source data entry had each year's
# Oct/Nov/Dec rows appended in the wrong place - they belong at the start
# of the block (the rows were off by one when the table was compiled).
# Fix: for every 12-row year block, rotate it so that the last three rows
# (Oct, Nov, Dec) move to the front, followed by the original first nine
# rows (Jan..Sep), for every column A:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$blockSize = 12
$numBlocks = 4
$lastCol = 8   # column H

for ($b = 0; $b -lt $numBlocks; $b++) {

    $blockStart = $firstDataRow + ($b * $blockSize)

    # Snapshot the 12 rows of this block (A:H) before overwriting anything.
    # NB: use Value2 (not Value) to read - Value2 round-trips the actual
    # cell content, Value does not.
    $snapshot = New-Object 'object[,]' $blockSize, $lastCol
    for ($r = 0; $r -lt $blockSize; $r++) {
        for ($c = 1; $c -le $lastCol; $c++) {
            $snapshot[$r, $c - 1] = $ws.Cells.Item($blockStart + $r, $c).Value2
        }
    }

    # Build the rotated order: old rows 9,10,11 (Oct,Nov,Dec) first,
    # then old rows 0..8 (Jan..Sep).
    $order = @(9, 10, 11, 0, 1, 2, 3, 4, 5, 6, 7, 8)

    for ($newR = 0; $newR -lt $blockSize; $newR++) {
        $oldR = $order[$newR]
        for ($c = 1; $c -le $lastCol; $c++) {
            $ws.Cells.Item($blockStart + $newR, $c).Value2 = $snapshot[$oldR, $c - 1]
        }
    }
}
